$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(2).Insert()

$ws.Range("A2").Value = "－"
$ws.Range("B2").Value = "2026年1月21日（令和8年1月21日）"
$ws.Range("C2").Value = "１令和８年度診療報酬改定に係る検討状況について`n２意見発表者による意見発表、中医協委員からの質問`n"
$ws.Range("D2").Value = "－"
$ws.Range("E2").Value = "資料`n"
$ws.Range("F2").Value = "－"

$ws.Range("C3").Value = "１令和８年度費用対効果評価制度の見直しについて`n２令和８年度保険医療材料制度の見直しについて`n３令和８年度薬価制度の見直しについて`n４市場拡大再算定について`n５歯科用貴金属価格の随時改定について`n６個別改定項目について（医療法等改正に伴う療養担当規則等の所要の見直しについて）`n７答申について（医療法等改正に伴う療養担当規則等の所要の見直しについて）`n８医療技術評価分科会からの報告について`n"
